# Update the cryptocurrency price/volume table on Sheet1 (cols B-E, rows 2-51)
# to match the latest scraped values, per the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'28.527.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.01%  "
# Row 3
$ws.Range("D3").Value = "'1.965.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.12%  "
# Row 4
$ws.Range("D4").Value = "'1.014"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.30%  "
# Row 5
$ws.Range("D5").Value = "'322.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.38%  "
# Row 6
$ws.Range("D6").Value = "'1.012"
$ws.Range("D6").Style = "Normal"
# Row 7
$ws.Range("D7").Value = "'0.4763"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.61%  "
# Row 8
$ws.Range("D8").Value = "'0.4048"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.80%  "
# Row 9
$ws.Range("D9").Value = "'54.09"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.34%  "
# Row 10
$ws.Range("D10").Value = "'0.08472"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.71%  "
# Row 11
$ws.Range("D11").Value = "'1.061"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.32%  "
# Row 12
$ws.Range("D12").Value = "'22.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.56%  "
# Row 13
$ws.Range("D13").Value = "'2.025.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.12%  "
# Row 14
$ws.Range("D14").Value = "'7.600"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.34%  "
# Row 15
$ws.Range("D15").Value = "'6.190"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.64%  "
# Row 16
$ws.Range("E16").Value = "  +0.18%  "
# Row 17
$ws.Range("D17").Value = "'90.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.32%  "
# Row 18
$ws.Range("D18").Value = "'0.00001069"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.84%  "
# Row 19
$ws.Range("D19").Value = "'0.06642"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.22%  "
# Row 20
$ws.Range("D20").Value = "'18.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.44%  "
# Row 21
$ws.Range("D21").Value = "'1.013"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.26%  "
# Row 22
$ws.Range("D22").Value = "'5.861"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.36%  "
# Row 23
$ws.Range("D23").Value = "'28.567.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.96%  "
# Row 24
$ws.Range("D24").Value = "'11.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.09%  "
# Row 25
$ws.Range("D25").Value = "'2.301"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "
# Row 26
$ws.Range("D26").Value = "'2.236.75"
$ws.Range("D26").Style = "Normal"
# Row 27
$ws.Range("D27").Value = "'155.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.57%  "
# Row 28
$ws.Range("D28").Value = "'20.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.22%  "
# Row 29
$ws.Range("D29").Value = "'5.895"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.98%  "
# Row 30
$ws.Range("D30").Value = "'2.166"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.96%  "
# Row 31
$ws.Range("D31").Value = "'124.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.72%  "
# Row 32
$ws.Range("D32").Value = "'0.9821"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.49%  "
# Row 33
$ws.Range("D33").Value = "'0.09642"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.06%  "
# Row 34
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.455"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.57%  "
# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.698"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.37%  "
# Row 36
$ws.Range("D36").Value = "'5.619"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.41%  "
# Row 37
$ws.Range("D37").Value = "'9.187"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.21%  "
# Row 38
$ws.Range("D38").Value = "'0.02333"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.89%  "
# Row 39
$ws.Range("D39").Value = "'0.06234"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.64%  "
# Row 40
$ws.Range("D40").Value = "'1.257"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.36%  "
# Row 41
$ws.Range("D41").Value = "'0.6210"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.46%  "
# Row 42
$ws.Range("D42").Value = "'11.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.39%  "
# Row 43
$ws.Range("D43").Value = "'1.012"
$ws.Range("D43").Style = "Normal"
# Row 44
$ws.Range("D44").Value = "'0.1914"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.25%  "
# Row 45
$ws.Range("D45").Value = "'1.352"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.74%  "
# Row 46
$ws.Range("D46").Value = "'0.5948"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.10%  "
# Row 47
$ws.Range("D47").Value = "'13.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.62%  "
# Row 48
$ws.Range("D48").Value = "'2.060"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.24%  "
# Row 49
$ws.Range("D49").Value = "'3.411"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.74%  "
# Row 50
$ws.Range("D50").Value = "'0.06821"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.67%  "
# Row 51
$ws.Range("D51").Value = "'111.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.16%  "
